$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the cDNA sample/protocol number from E7420 to E7420L.
#    This value is shared by every row in column G (G2:G41).
$ws.Range("G2:G41").Value = "E7420L"

# 2. Replace the "=FALSE()" formulas in the roboticS2Prep column (H2:H41)
#    with a plain boolean FALSE literal (no formula), as part of the
#    accuracy check update.
$ws.Range("H2:H41").Value = $false

# 3. Scroll the sheet view up one row so row 19 becomes the top-left
#    visible cell (was row 20).
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
